# Auto-generated edit script: refresh cached FFXIV market-price derived
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) across all 8
# class/job sheets, per the scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2030.2413
$ws.Range("J112").Value = 2065.8147
$ws.Range("L112").Value = 6197.4441
$ws.Range("N112").Value = -8413.444100000001
$ws.Range("H116").Value = 5417.5
$ws.Range("I116").Value = 4199.6665
$ws.Range("J116").Value = 6635.3335
$ws.Range("K116").Value = 4199.6665
$ws.Range("L116").Value = 6635.3335
$ws.Range("M116").Value = -757.6665000000003
$ws.Range("N116").Value = -13519.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 26792.176
$ws.Range("I45").Value = 29562.666
$ws.Range("J45").Value = 6013.5
$ws.Range("K45").Value = 29562.666
$ws.Range("L45").Value = 6013.5
$ws.Range("M45").Value = -29185.666
$ws.Range("N45").Value = -6767.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3017.4443
$ws.Range("I86").Value = 2679.5715
$ws.Range("J86").Value = 4200
$ws.Range("K86").Value = 2679.5715
$ws.Range("L86").Value = 4200
$ws.Range("M86").Value = -1556.5715
$ws.Range("N86").Value = -6446
$ws.Range("H89").Value = 3017.4443
$ws.Range("I89").Value = 2679.5715
$ws.Range("J89").Value = 4200
$ws.Range("K89").Value = 13397.8575
$ws.Range("L89").Value = 21000
$ws.Range("M89").Value = -7781.8575
$ws.Range("N89").Value = -32232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1423.4828
$ws.Range("I16").Value = 1232.4348
$ws.Range("J16").Value = 2155.8333
$ws.Range("K16").Value = 1232.4348
$ws.Range("L16").Value = 2155.8333
$ws.Range("M16").Value = -945.4348
$ws.Range("N16").Value = -2729.8333
$ws.Range("H22").Value = 2297
$ws.Range("I22").Value = 2297
$ws.Range("K22").Value = 2297
$ws.Range("M22").Value = -1947
$ws.Range("H86").Value = 4127
$ws.Range("I86").Value = 4166.6665
$ws.Range("J86").Value = 4008
$ws.Range("K86").Value = 4166.6665
$ws.Range("L86").Value = 4008
$ws.Range("M86").Value = -3043.6665
$ws.Range("N86").Value = -6254
$ws.Range("H89").Value = 4127
$ws.Range("I89").Value = 4166.6665
$ws.Range("J89").Value = 4008
$ws.Range("K89").Value = 20833.3325
$ws.Range("L89").Value = 20040
$ws.Range("M89").Value = -15217.3325
$ws.Range("N89").Value = -31272
$ws.Range("H113").Value = 1423.4828
$ws.Range("I113").Value = 1232.4348
$ws.Range("J113").Value = 2155.8333
$ws.Range("K113").Value = 1232.4348
$ws.Range("L113").Value = 2155.8333
$ws.Range("M113").Value = 937.5652
$ws.Range("N113").Value = -6495.8333
$ws.Range("H132").Value = 2351.6428
$ws.Range("I132").Value = 1658.7778
$ws.Range("J132").Value = 3598.8
$ws.Range("K132").Value = 4976.3334
$ws.Range("L132").Value = 10796.4
$ws.Range("M132").Value = -2446.3334
$ws.Range("N132").Value = -15856.4
$ws.Range("H133").Value = 71997.336
$ws.Range("J133").Value = 71997.336
$ws.Range("L133").Value = 71997.336
$ws.Range("N133").Value = -77057.336
$ws.Range("H134").Value = 3189.0588
$ws.Range("I134").Value = 3080.9333
$ws.Range("K134").Value = 9242.7999
$ws.Range("M134").Value = -6707.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 5330.6665
$ws.Range("I22").Value = 4498
$ws.Range("J22").Value = 5747
$ws.Range("K22").Value = 13494
$ws.Range("L22").Value = 17241
$ws.Range("M22").Value = -13325
$ws.Range("N22").Value = -17579
$ws.Range("H27").Value = 5330.6665
$ws.Range("I27").Value = 4498
$ws.Range("J27").Value = 5747
$ws.Range("K27").Value = 13494
$ws.Range("L27").Value = 17241
$ws.Range("M27").Value = -13392
$ws.Range("N27").Value = -17445
$ws.Range("H38").Value = 524.0909
$ws.Range("I38").Value = 148.33333
$ws.Range("K38").Value = 444.99999
$ws.Range("M38").Value = -97.99998999999997
$ws.Range("H107").Value = 613.4
$ws.Range("J107").Value = 613.4
$ws.Range("L107").Value = 1840.2
$ws.Range("N107").Value = -5680.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3544.5334
$ws.Range("I122").Value = 2262.2273
$ws.Range("K122").Value = 6786.6819
$ws.Range("M122").Value = -4336.6819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2643.95
$ws.Range("I7").Value = 3186.3333
$ws.Range("J7").Value = 1830.375
$ws.Range("K7").Value = 3186.3333
$ws.Range("L7").Value = 1830.375
$ws.Range("M7").Value = -3074.3333
$ws.Range("N7").Value = -2054.375
$ws.Range("H46").Value = 2666.6667
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376
$ws.Range("H126").Value = 2643.95
$ws.Range("I126").Value = 3186.3333
$ws.Range("J126").Value = 1830.375
$ws.Range("K126").Value = 9558.999899999999
$ws.Range("L126").Value = 5491.125
$ws.Range("M126").Value = -7088.999899999999
$ws.Range("N126").Value = -10431.125
$ws.Range("H134").Value = 103991.336
$ws.Range("J134").Value = 103991.336
$ws.Range("L134").Value = 103991.336
$ws.Range("N134").Value = -114131.336
$ws.Range("H135").Value = 54198.8
$ws.Range("J135").Value = 54198.8
$ws.Range("L135").Value = 54198.8
$ws.Range("N135").Value = -64338.8
$ws.Range("H137").Value = 57630.895
$ws.Range("J137").Value = 57630.895
$ws.Range("L137").Value = 57630.895
$ws.Range("N137").Value = -67830.89499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8151
$ws.Range("I62").Value = 3877
$ws.Range("J62").Value = 9219.5
$ws.Range("K62").Value = 3877
$ws.Range("L62").Value = 9219.5
$ws.Range("M62").Value = -3253
$ws.Range("N62").Value = -10467.5
$ws.Range("H65").Value = 8151
$ws.Range("I65").Value = 3877
$ws.Range("J65").Value = 9219.5
$ws.Range("K65").Value = 19385
$ws.Range("L65").Value = 46097.5
$ws.Range("M65").Value = -16265
$ws.Range("N65").Value = -52337.5
$ws.Range("H96").Value = 4700
$ws.Range("J96").Value = 5999.3335
$ws.Range("L96").Value = 5999.3335
$ws.Range("N96").Value = -8745.333500000001
$ws.Range("H107").Value = 769
$ws.Range("I107").Value = 640.2
$ws.Range("K107").Value = 1920.6
$ws.Range("M107").Value = -0.6000000000001364
$ws.Range("H113").Value = 460.3846
$ws.Range("I113").Value = 272.375
$ws.Range("J113").Value = 761.2
$ws.Range("K113").Value = 817.125
$ws.Range("L113").Value = 2283.6
$ws.Range("M113").Value = 1352.875
$ws.Range("N113").Value = -6623.6
$ws.Range("H122").Value = 22729126
$ws.Range("I122").Value = 1674.375
$ws.Range("K122").Value = 5023.125
$ws.Range("M122").Value = -2573.125
$ws.Range("H126").Value = 2980.4285
$ws.Range("I126").Value = 1753.5
$ws.Range("J126").Value = 4616.3335
$ws.Range("K126").Value = 5260.5
$ws.Range("L126").Value = 13849.0005
$ws.Range("M126").Value = -2790.5
$ws.Range("N126").Value = -18789.0005
$ws.Range("H132").Value = 2305.8948
$ws.Range("I132").Value = 2118.4285
$ws.Range("J132").Value = 4493
$ws.Range("K132").Value = 6355.2855
$ws.Range("L132").Value = 13479
$ws.Range("M132").Value = -3825.2855
$ws.Range("N132").Value = -18539
$ws.Range("H136").Value = 66670550
$ws.Range("J136").Value = 13125
$ws.Range("L136").Value = 39375
$ws.Range("N136").Value = -44475
